# Weekly price-sheet update: a new daily record is published for
# "Terminal Hortofrutícola Agro Chillán" (Frutilla), which pushes a brand
# new row in right after the existing row 52 (i.e. at row 53), shifting
# every subsequent record down by one row and extending the used range
# from A1:T121 to A1:T122.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 53; rows 53:121 shift down to 54:122.
$ws.Rows("53:53").Insert()

# Fill the newly inserted row with the new weekly record.
$ws.Range("A53").Value = 7
$ws.Range("B53").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C53").Value = "Ñuble"
$ws.Range("D53").Value = "2021-10-08"
$ws.Range("E53").Value = 16
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100101
$ws.Range("H53").Value = "Berries"
$ws.Range("I53").Value = 100112025
$ws.Range("J53").Value = "Frutilla"
$ws.Range("K53").Value = "Sin especificar"
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 60
$ws.Range("N53").Value = 12000
$ws.Range("O53").Value = 13000
$ws.Range("P53").Value = 12500
$ws.Range("Q53").Value = "`$/bandeja 7 kilos"
$ws.Range("R53").Value = "Provincia de Melipilla"
$ws.Range("S53").Value = 1786
$ws.Range("T53").Value = 7
